# 12-5-23 late updates with Pete
#
# The "micron"/"um" unit rows get moved up (right after "microinch") in the
# unit list, and every factor in column B that used to be a live formula
# (=1/x) is flattened down to its last computed (static) value. Finally the
# sheet's last active selection moves from A21 to B20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-label rows 15-18 (A column) so "micron" / "um" sit right after
#     "microinch", pushing "mi (US survey)" / "ft (US survey)" down. The
#     B-column factors stay fixed to their row and just lose their formula,
#     becoming the plain static number that used to live there.
$ws.Range("A15").Value = "micron"
$ws.Range("A16").Value = "&#181m"
$ws.Range("A17").Value = "mi (US survey)"
$ws.Range("A18").Value = "ft (US survey)"

# --- Flatten every remaining "=1/x" formula in column B to its static
#     cached value (rows 8-11, 13-18). Rows 15-18 also pick up the values
#     that correspond to their re-ordered labels above.
$ws.Range("B8").Value = 3.280833436679587
$ws.Range("B9").Value = 39.370078740157481
$ws.Range("B10").Value = 0.00062137119223733392
$ws.Range("B11").Value = 1.0936132983377078
$ws.Range("B13").Value = 39370.078740157478
$ws.Range("B14").Value = 39370078.740157485
$ws.Range("B15").Value = 1000000
$ws.Range("B16").Value = 1000000
$ws.Range("B17").Value = 0.00062137003393301752
$ws.Range("B18").Value = 3.280833436679587

# --- Move the sheet's remembered selection from A21 to B20.
$ws.Range("B20").Select() | Out-Null
